$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets carry an identical copy of this table,
# so the same update needs to be applied to each of them.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Column B holds the start date as literal text (e.g. "2024.02.07").
    # Force a Text number format first so Excel doesn't auto-convert the
    # new "2024-02-07" style strings into real date serials.
    $ws.Range("B2:B5").NumberFormat = "@"

    # Row 2: date separator dot -> dash, swap in the "LPJ" event name,
    # bump the interest count, clear the old bilibili link.
    $ws.Range("B2").Value = "2024-02-07"
    $ws.Range("C2").Value = "丽水·LPJ 现实X次元动漫展"
    $ws.Range("F2").Value = 271
    $ws.Range("I2").Value = ""

    # Row 3: date separator dot -> dash, swap in the "YA" event name,
    # bump the interest count, clear the old bilibili link.
    $ws.Range("B3").Value = "2024-02-07"
    $ws.Range("C3").Value = "丽水·YA●怀旧only"
    $ws.Range("F3").Value = 221
    $ws.Range("I3").Value = ""

    # Row 4: date separator dot -> dash, clear the event name,
    # bump the interest count, clear the old bilibili link.
    $ws.Range("B4").Value = "2024-02-14"
    $ws.Range("C4").Value = ""
    $ws.Range("F4").Value = 17
    $ws.Range("I4").Value = ""

    # Row 5: date separator dot -> dash, swap in the "崩铁" event name,
    # bump the interest count, clear the old bilibili link.
    $ws.Range("B5").Value = "2024-02-18"
    $ws.Range("C5").Value = "龙泉·崩X铁X原ONLY"
    $ws.Range("F5").Value = 259
    $ws.Range("I5").Value = ""
}
